$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

Set-TextValue $ws 'D2' '246.84'
Set-TextValue $ws 'E2' '0.84%'
Set-TextValue $ws 'D3' '29.54'
Set-TextValue $ws 'E3' '7.67%'
Set-TextValue $ws 'D4' '5.194'
Set-TextValue $ws 'E4' '1.54%'
Set-TextValue $ws 'E5' '0.85%'
Set-TextValue $ws 'D6' '6.562'
Set-TextValue $ws 'E6' '0.62%'
Set-TextValue $ws 'D7' '0.8589'
Set-TextValue $ws 'E7' '4.76%'
Set-TextValue $ws 'D8' '0.8665'
Set-TextValue $ws 'E8' '1.67%'
Set-TextValue $ws 'E9' '2.20%'
Set-TextValue $ws 'D10' '0.07090'
Set-TextValue $ws 'E10' '2.01%'
Set-TextValue $ws 'D11' '0.03015'
Set-TextValue $ws 'E11' '4.69%'
Set-TextValue $ws 'D12' '0.09380'
Set-TextValue $ws 'E12' '-0.13%'
Set-TextValue $ws 'D13' '0.001526'
Set-TextValue $ws 'E13' '0.74%'
Set-TextValue $ws 'D14' '0.0005979'
Set-TextValue $ws 'E14' '-94.12%'
Set-TextValue $ws 'D15' '0.006031'
Set-TextValue $ws 'E15' '-2.95%'
Set-TextValue $ws 'E16' '5,225.03%'
Set-TextValue $ws 'D17' '3.493'
Set-TextValue $ws 'E17' '-0.56%'
Set-TextValue $ws 'D18' '3.100'
Set-TextValue $ws 'E18' '3.01%'
Set-TextValue $ws 'D19' '2.283'
Set-TextValue $ws 'E19' '-1.53%'
Set-TextValue $ws 'E20' '-0.16%'
Set-TextValue $ws 'D21' '0.03317'
Set-TextValue $ws 'E21' '2.95%'
Set-TextValue $ws 'E22' '1.21%'
Set-TextValue $ws 'D23' '3.473'
Set-TextValue $ws 'E23' '-2.39%'
Set-TextValue $ws 'D24' '0.04144'
Set-TextValue $ws 'E24' '2.88%'
Set-TextValue $ws 'E25' '0.44%'
Set-TextValue $ws 'D26' '0.001227'
Set-TextValue $ws 'E26' '1.01%'
Set-TextValue $ws 'D27' '0.004997'
Set-TextValue $ws 'E27' '11.57%'
Set-TextValue $ws 'D28' '0.0001211'
Set-TextValue $ws 'E28' '2.59%'
Set-TextValue $ws 'D40' '0.03753'
Set-TextValue $ws 'E40' '0.94%'
Set-TextValue $ws 'B41' 'KickToken'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws 'D41' '0.005813'
Set-TextValue $ws 'E41' '-2.96%'
Set-TextValue $ws 'B42' 'BKEXToken'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws 'D42' '0.1071'
Set-TextValue $ws 'E42' '1.28%'
Set-TextValue $ws 'B43' 'CEJI'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws 'D43' '0.002101'
Set-TextValue $ws 'E43' '-10.59%'
Set-TextValue $ws 'D44' '0.009497'
Set-TextValue $ws 'E44' '-2.26%'
Set-TextValue $ws 'D45' '0.00005288'
Set-TextValue $ws 'E45' '3.55%'
Set-TextValue $ws 'E46' '0.06%'
Set-TextValue $ws 'D47' '0.05699'
Set-TextValue $ws 'E47' '-43.58%'
Set-TextValue $ws 'D48' '0.002279'
Set-TextValue $ws 'E48' '-9.43%'
Set-TextValue $ws 'D49' '0.00002101'
Set-TextValue $ws 'E49' '0.06%'
Set-TextValue $ws 'D50' '0.0002001'
Set-TextValue $ws 'E50' '0.06%'
